$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 10 data
$ws.Range("A10").Value = "Employment by industry"
$ws.Range("B10").Value = "<a href='https://www.nomisweb.co.uk/datasets/apsnew'>Annual Population Survey</a>"
$ws.Range("C10").Value = "Jul 2021 - Jun 2022 (11/10/22)"
$ws.Range("D10").Value = "Oct 2021 - Sep 2022 (17/01/23)"

# Match the style of column B used in other rows (link style) by copying formats only
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update sheet view / selection: remove frozen/scrolled topLeftCell and change selection to C12
$wb.Windows.Item(1).ScrollColumn = 1
$ws.Range("C12").Select()
